# Refatorando o consolidador para modelo ETL
# Atualiza os dados de absenteismo (linhas 2-11) conforme o novo dataset

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 81069; B = "André Novaes";          C = "Recursos Humanos";        D = "Viagem de negócios"; E = 5; F = 45078;  G = 12303.11 },
    @{ Row = 3;  A = 66874; B = "Alana Caldeira";         C = "Engenharia";              D = "Problemas pessoais"; E = 7; F = 45098;  G = 2957 },
    @{ Row = 4;  A = 46512; B = "Pedro Rocha";            C = "Operações";               D = "Viagem de negócios"; E = 6; F = 45099;  G = 6197.14 },
    @{ Row = 5;  A = 46792; B = "Francisco Fernandes";    C = "Marketing";               D = "Doença";             E = 1; F = 45100;  G = 3926.29 },
    @{ Row = 6;  A = 1961;  B = "Marcelo Nunes";          C = "Atendimento ao Cliente";  D = "Outros";             E = 5; F = 45087;  G = 11525.04 },
    @{ Row = 7;  A = 66415; B = "Rafaela Ribeiro";        C = "Marketing";               D = "Doença";             E = 4; F = 45094;  G = 8298.940000000001 },
    @{ Row = 8;  A = 46209; B = "Mirella Cavalcanti";     C = "Atendimento ao Cliente";  D = "Doença";             E = 7; F = 45102;  G = 7227.92 },
    @{ Row = 9;  A = 1980;  B = "Melissa da Paz";         C = "Financeiro";              D = "Problemas pessoais"; E = 2; F = 45096;  G = 3914.37 },
    @{ Row = 10; A = 88259; B = "Alícia Barbosa";         C = "Recursos Humanos";        D = "Viagem de negócios"; E = 1; F = 45080;  G = 9903.57 },
    @{ Row = 11; A = 9586;  B = "Maria Vitória Martins";  C = "Marketing";               D = "Problemas pessoais"; E = 4; F = 45084;  G = 10735.67 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
